$wb = $excel.ActiveWorkbook

# --- "Customer Issue Log" becomes the Data Table used as a Queue: rebuild
#     its header row with the new Email-centric columns (Subject/Body split
#     out, Status renamed/simplified) ---
$ws1 = $wb.Worksheets.Item("Customer Issue Log")

$ws1.Range("A1").Value = "Date"
$ws1.Range("D1").Value = "Email Body"
$ws1.Range("B1").Value = "Customer Email"
$ws1.Range("F1").Value = "Status"
$ws1.Range("C1").Value = "Email Subject"
$ws1.Range("E1").Value = "AI Genereated Response"
$ws1.Range("G1").Value = "Issue Category"
$ws1.Range("H1").Value = "Date of Resolution"

# Carry the bold header style onto the two newly-added columns (G:H) by
# copying the format from an already-styled header cell.
$ws1.Range("A1").Copy() | Out-Null
$ws1.Range("G1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws1.Range("C1:H1").EntireColumn.AutoFit() | Out-Null

# New selection/active-cell on this sheet, and make it the active tab.
$ws1.Range("M5").Select() | Out-Null
$ws1.Activate() | Out-Null
